$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1666666666666667
$ws.Range("C2").Value = 0.5
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.1666666666666667
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("B6").Value = 0.03225806451612903
$ws.Range("F6").Value = 0.0967741935483871
$ws.Range("J6").Value = 0.1612903225806452
$ws.Range("Q6").Value = 0.0967741935483871
$ws.Range("R6").Value = 0.06451612903225806
$ws.Range("S6").Value = 0.5483870967741935
$ws.Range("B7").Value = 0.0625
$ws.Range("F7").Value = 0.09375
$ws.Range("J7").Value = 0.125
$ws.Range("Q7").Value = 0.0625
$ws.Range("R7").Value = 0.09375
$ws.Range("S7").Value = 0.5625
$ws.Range("B8").Value = 0.0303030303030303
$ws.Range("F8").Value = 0.07575757575757576
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("Q8").Value = 0.09090909090909091
$ws.Range("R8").Value = 0.07575757575757576
$ws.Range("S8").Value = 0.6363636363636364
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.04
$ws.Range("Q9").Value = 0.08
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.64
$ws.Range("B10").Value = 0.04210526315789474
$ws.Range("D10").Value = 0.02105263157894737
$ws.Range("F10").Value = 0.1157894736842105
$ws.Range("J10").Value = 0.07368421052631578
$ws.Range("O10").Value = 0.02105263157894737
$ws.Range("Q10").Value = 0.07368421052631578
$ws.Range("R10").Value = 0.1263157894736842
$ws.Range("S10").Value = 0.5263157894736842
$ws.Range("G11").Value = 0.1363636363636364
$ws.Range("J11").Value = 0.06818181818181818
$ws.Range("K11").Value = 0.2045454545454546
$ws.Range("L11").Value = 0.5909090909090909
$ws.Range("G12").Value = 0.8846153846153846
$ws.Range("J12").Value = 0.07692307692307693
$ws.Range("S12").Value = 0.03846153846153846
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.125
$ws.Range("H15").Value = 0.2941176470588235
$ws.Range("I15").Value = 0.1764705882352941
$ws.Range("J15").Value = 0.05882352941176471
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.3529411764705883
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("J16").Value = 0.5
$ws.Range("S16").Value = 0.3333333333333333
$ws.Range("H17").Value = 0.2857142857142857
$ws.Range("I17").Value = 0.09523809523809523
$ws.Range("J17").Value = 0.2857142857142857
$ws.Range("K17").Value = 0.1904761904761905
$ws.Range("M17").Value = 0.04761904761904762
$ws.Range("S17").Value = 0.09523809523809523
$ws.Range("F18").Value = 0.04166666666666666
$ws.Range("H18").Value = 0.375
$ws.Range("I18").Value = 0.125
$ws.Range("K18").Value = 0.1666666666666667
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01612903225806452
$ws.Range("H19").Value = 0.2473118279569892
$ws.Range("I19").Value = 0.0913978494623656
$ws.Range("J19").Value = 0.2688172043010753
$ws.Range("K19").Value = 0.1397849462365591
$ws.Range("M19").Value = 0.03763440860215054
$ws.Range("O19").Value = 0.05376344086021505
$ws.Range("S19").Value = 0.1451612903225807
